# Reworked initial set of functions
# - Drop the "value.alternative" column (E) from the "vars" and "pars" sheets
# - Rename the "value.default" header (D1) to "default" on those sheets
# - Make "pars" the active sheet/tab; "eqns" is no longer the active sheet
# - Update remembered selections on the "pars" and "eqns" sheets

$wb = $excel.ActiveWorkbook

$wsVars = $wb.Worksheets.Item("vars")
$wsVars.Columns.Item(5).Delete()
$wsVars.Range("D1").Value = "default"

$wsPars = $wb.Worksheets.Item("pars")
$wsPars.Columns.Item(5).Delete()
$wsPars.Range("D1").Value = "default"

$wsEqns = $wb.Worksheets.Item("eqns")
[void]$wsEqns.Range("F6").Select()

$wsPars.Activate()
[void]$wsPars.Range("D2").Select()
